$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# sheet "species"
# ---------------------------------------------------------------------------
$species = $wb.Worksheets.Item("species")

# numeric value tweaks (plain cells, no explicit style -> safe with .Value)
$species.Range("F3").Value = 1
$species.Range("D4").Value = 0
$species.Range("D5").Value = 0

# B5/C5 species name stays "C" (text unchanged, shared-string table just gets
# recompacted once the reactions-sheet strings are rewritten below)
$species.Range("B5").Value = "C"
$species.Range("C5").Value = "C"

# ---------------------------------------------------------------------------
# sheet "reactions"
# ---------------------------------------------------------------------------
$reactions = $wb.Worksheets.Item("reactions")

# Remove the blank formatting row 6 - everything below (rows 7..35) shifts up
# by one (row 35 disappears, dimension shrinks from L35 to L34).
$reactions.Rows(6).Delete()

# Row 4: rule text "A ===> B" -> "A => B", weight 0.25 -> 1
# (use a scratch cell + PasteSpecial(values) so the destination's existing
# cell style, s="11"/quotePrefix, survives the write - a plain .Value=
# assignment resets the xf index on this engine).
$reactions.Range("ZZ1").Value = "A => B"
$reactions.Range("ZZ1").Copy()
$reactions.Range("C4").PasteSpecial(-4163)
$reactions.Range("ZZ1").Clear()
$reactions.Range("D4").Value = 1

# Row 5: rule text "!B => C" -> "B ===> C", weight 0.25 -> 1, and the B5 label
# goes back to "r2" (same text as before, index just gets recompacted)
$reactions.Range("B5").Value = "r2"
$reactions.Range("ZZ1").Value = "B ===> C"
$reactions.Range("ZZ1").Copy()
$reactions.Range("C5").PasteSpecial(-4163)
$reactions.Range("ZZ1").Clear()
$reactions.Range("D5").Value = 1

# ---------------------------------------------------------------------------
# view / selection state
# ---------------------------------------------------------------------------
$species.Range("F5").Select()
$reactions.Select()
$reactions.Range("B5").Select()
